# Swap the bond detail columns (shortName, marketValue, couponRate, faceAmount,
# isin, percentWeight, cusip, sedol) between paired rows so that each pair's
# two holdings line up correctly with their respective bond identifiers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "F", "G", "I", "K", "L", "Q", "R")
$pairs = @(
    @(37, 38),
    @(41, 42),
    @(45, 46),
    @(47, 48),
    @(51, 52),
    @(55, 56)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")
        $v1 = $cell1.Text
        $v2 = $cell2.Text

        # Numeric-looking values (e.g. "2.875", "173101346") would otherwise be
        # silently re-typed as real numbers when assigned back through COM,
        # which would lose the original text formatting (trailing zeros, etc.)
        # and cell data type. Prefixing with an apostrophe forces Excel to
        # store them as literal text, matching the source workbook where
        # every cell in this sheet is an inline/shared string.
        $isNum1 = $v1 -match '^-?[0-9]+(\.[0-9]+)?$'
        $isNum2 = $v2 -match '^-?[0-9]+(\.[0-9]+)?$'

        if ($isNum1) { $cell2.Value2 = "'" + $v1 } else { $cell2.Value2 = $v1 }
        if ($isNum2) { $cell1.Value2 = "'" + $v2 } else { $cell1.Value2 = $v2 }
    }
}
